$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bosquets")

# Delete column B ("Arreglo aleatorio óptimo bosquets") which shifts remaining columns left
$ws.Range("B1").EntireColumn.Delete()

# Update header row (now A:E after the shift)
$ws.Range("A1").Value = "Estimador óptimo"
$ws.Range("B1").Value = "MAE bosquets"
$ws.Range("C1").Value = "MSE bosquets"
$ws.Range("D1").Value = "RMSE bosquets"
$ws.Range("E1").Value = "R2 bosquets"

# Update data row 2 with new values
$ws.Range("A2").Value = 102
$ws.Range("B2").Value = 0.6082810349699448
$ws.Range("C2").Value = 0.5844243957844734
$ws.Range("D2").Value = 0.7644765501861214
$ws.Range("E2").Value = 0.03229832537117805
